# Changed ui layout, moved update logic to app
#
# The "Translation" worksheet's GB-* columns (G/H/I) are collapsed into a
# single "GB" follow-up column (G), several rows are re-pointed at
# different (and new) string/typography ids, and the three trailing rows
# that only existed to host the old GB-ALIGNMENT / GB-DIRECTION data are
# cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Row 3 (header row): G3 keeps "GB-TYPOGRAPHY", H3/I3 are dropped ---
$ws.Range("G3").Value = "GB-TYPOGRAPHY"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

# --- Row 4 ---
$ws.Range("F4").Value = "rpm"

# --- Row 5 ---
$ws.Range("B5").Value = "SingleUseId4"
$ws.Range("C5").Value = "Default"
$ws.Range("F5").Value = "km"

# --- Row 6 ---
$ws.Range("B6").Value = "currentRpm"
$ws.Range("D6").Value = "Right"
$ws.Range("F6").Value = "<number>"
$ws.Range("G6").Value = "Big"

# --- Row 7 ---
$ws.Range("B7").Value = "totalDistance"
$ws.Range("G7").Value = "Large"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""

# --- Row 8 ---
$ws.Range("B8").Value = "SingleUseId6"
$ws.Range("C8").Value = "Big"
# F8 needs to literally be the text "0" (a format-string id), not the
# number 0 - go through a TEXT() helper cell + paste-values so Excel
# keeps it as a shared string without bolting a new text number-format
# onto the cell's style.
$ws.Range("ZZ1").Formula = '=TEXT(0,"0")'
$ws.Range("ZZ1").Copy()
$ws.Range("F8").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("G8").Value = "Big"
$ws.Range("H8").Value = ""

# --- Row 9 ---
$ws.Range("B9").Value = "SingleUseId7"
$ws.Range("C9").Value = "Large"
$ws.Range("D9").Value = "Left"
# Same story as F8 above, but the literal text is "000000".
$ws.Range("ZZ1").Formula = '=TEXT(0,"000000")'
$ws.Range("ZZ1").Copy()
$ws.Range("F9").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$ws.Range("H9").Value = ""

$excel.CutCopyMode = $false

# --- Rows 10-12 no longer hold any data ---
$ws.Range("B10:H10").Value = ""
$ws.Range("B11:H11").Value = ""
$ws.Range("B12:H12").Value = ""
